$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.09"
$ws.Range("E2").Value = "'0.82%"
$ws.Range("D3").Value = "'27.29"
$ws.Range("E3").Value = "'0.43%"
$ws.Range("D4").Value = "'4.841"
$ws.Range("E4").Value = "'1.76%"
$ws.Range("D5").Value = "'0.06322"
$ws.Range("E5").Value = "'0.30%"
$ws.Range("E6").Value = "'1.42%"
$ws.Range("D7").Value = "'1.326"
$ws.Range("E7").Value = "'-2.58%"
$ws.Range("D8").Value = "'0.8883"
$ws.Range("E8").Value = "'1.33%"
$ws.Range("D9").Value = "'0.1510"
$ws.Range("E9").Value = "'-0.04%"
$ws.Range("D10").Value = "'0.05384"
$ws.Range("E10").Value = "'6.19%"
$ws.Range("D11").Value = "'0.07442"
$ws.Range("E11").Value = "'-2.48%"
$ws.Range("D12").Value = "'0.02899"
$ws.Range("E12").Value = "'-2.37%"
$ws.Range("E13").Value = "'-0.54%"
$ws.Range("D14").Value = "'0.001562"
$ws.Range("E14").Value = "'0.07%"
$ws.Range("D15").Value = "'0.0006349"
$ws.Range("E15").Value = "'-0.23%"
$ws.Range("D16").Value = "'0.006031"
$ws.Range("E16").Value = "'0.90%"
$ws.Range("D17").Value = "'3.473"
$ws.Range("E17").Value = "'0.80%"
$ws.Range("D18").Value = "'3.297"
$ws.Range("E18").Value = "'-0.02%"
$ws.Range("D19").Value = "'2.234"
$ws.Range("E19").Value = "'-1.69%"
$ws.Range("E21").Value = "'1.24%"
$ws.Range("D22").Value = "'3.916"
$ws.Range("E22").Value = "'0.14%"
$ws.Range("D23").Value = "'0.1506"
$ws.Range("E23").Value = "'9.14%"
$ws.Range("D24").Value = "'0.04383"
$ws.Range("E24").Value = "'-0.10%"
$ws.Range("E25").Value = "'0.18%"
$ws.Range("D26").Value = "'0.004251"
$ws.Range("E26").Value = "'10.71%"
$ws.Range("E28").Value = "'-1.67%"
$ws.Range("E29").Value = "'-14.81%"
$ws.Range("D40").Value = "'0.03995"
$ws.Range("E40").Value = "'-2.48%"
$ws.Range("D41").Value = "'0.006651"
$ws.Range("E41").Value = "'-2.06%"
$ws.Range("E42").Value = "'20.31%"
$ws.Range("D43").Value = "'0.002150"
$ws.Range("E43").Value = "'3.36%"
$ws.Range("E44").Value = "'2.39%"
$ws.Range("D45").Value = "'0.00005528"
$ws.Range("E45").Value = "'6.77%"
$ws.Range("E46").Value = "'9.53%"
$ws.Range("E47").Value = "'-19.66%"
